$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date placeholder text from
#    14/09/2022 to 23/09/2023 on the slide master and every slide layout.
# ---------------------------------------------------------------------------
$oldDate = "14/09/2022"
$newDate = "23/09/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $isDatePh = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
            if ($isDatePh -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Bump the referenced PHP version on slide 1 from 8.1 to 8.2 wherever it
#    is mentioned ("Apache webserver + PHP 8.1" -> "Apache webserver + PHP 8.2").
# ---------------------------------------------------------------------------
$oldPhp = "Apache webserver + PHP 8.1"
$newPhp = "Apache webserver + PHP 8.2"

$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldPhp) {
            $sh.TextFrame.TextRange.Text = $newPhp
        }
    }
}
